# Updated cryptos list on Sat May 18 09:51:23 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures on the
# active worksheet to the latest scraped coinranking.com values.
#
# For cells whose new text happens to look like a plain number (e.g. 
# "1.00", "6.50"), NumberFormat is temporarily forced to Text so Excel
# keeps the exact string (incl. trailing zeros) instead of coercing it to
# a float; the style is then reset back to Normal so no stray number
# format is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.095.99'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '3.117.69'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.114.40'
$ws.Range("E8").Value = '  +2.83%  '
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.50'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").Value = '3.632.49'
$ws.Range("E16").Value = '  +2.80%  '
$ws.Range("D17").Value = '67.137.64'
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("D19").Value = '3.119.06'
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '487.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.717'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("E25").Value = '  +3.86%  '
$ws.Range("E26").Value = '  +3.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("E31").Value = '  +1.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.28%  '
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  +0.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.990'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("E39").Value = '  +3.36%  '
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("E41").Value = '  +1.62%  '
$ws.Range("E42").Value = '  +1.91%  '
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("D45").Value = '2.846.43'
$ws.Range("E45").Value = '  +4.57%  '
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '383.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '137.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  +2.60%  '
$ws.Range("E51").Value = '  +0.07%  '
